$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 668
$ws1.Range("F3").Value = 508
$ws1.Range("F4").Value = 37
$ws1.Range("F7").Value = 44
$ws1.Range("F8").Value = 2660
$ws1.Range("F9").Value = 4174
$ws1.Range("F10").Value = 102

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 61

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 668
$ws4.Range("F3").Value = 508
$ws4.Range("F4").Value = 37
$ws4.Range("F7").Value = 44
$ws4.Range("F8").Value = 2660
$ws4.Range("F9").Value = 4174
$ws4.Range("F10").Value = 102
$ws4.Range("F11").Value = 61

$wb.Save()
